$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 3.704051282521292
$ws.Range("C2").Value = 2.350769727091361
$ws.Range("D2").Value = 2.608743059361374
$ws.Range("E2").Value = 1.467857330083988
$ws.Range("F2").Value = 3.190166447130025
$ws.Range("G2").Value = 3.360965967234947
$ws.Range("H2").Value = 3.333022892298317
$ws.Range("B3").Value = -1.35328155542993
$ws.Range("C3").Value = -1.095308223159918
$ws.Range("D3").Value = -2.236193952437304
$ws.Range("E3").Value = -0.5138848353912668
$ws.Range("F3").Value = -0.3430853152863449
$ws.Range("G3").Value = -0.371028390222975
$ws.Range("B4").Value = 0.2579733322700122
$ws.Range("C4").Value = -0.8829123970073738
$ws.Range("D4").Value = 0.8393967200386636
$ws.Range("E4").Value = 1.010196240143586
$ws.Range("F4").Value = 0.9822531652069555
$ws.Range("G4").Value = -0.0757466556224102
$ws.Range("H4").Value = 1.072433058307495
$ws.Range("I4").Value = 0.4024492143100105
$ws.Range("J4").Value = 0.1212484254788393
$ws.Range("B5").Value = -1.140885729277386
$ws.Range("C5").Value = 0.5814233877686514
$ws.Range("D5").Value = 0.7522229078735734
$ws.Range("E5").Value = 0.7242798329369433
$ws.Range("F5").Value = -0.3337199878924224
$ws.Range("G5").Value = 0.8144597260374833
$ws.Range("H5").Value = 0.1444758820399983
$ws.Range("I5").Value = -0.1367249067911729
$ws.Range("B6").Value = 1.722309117046037
$ws.Range("C6").Value = 1.893108637150959
$ws.Range("D6").Value = 1.865165562214329
$ws.Range("E6").Value = 0.8071657413849636
$ws.Range("F6").Value = 1.955345455314869
$ws.Range("G6").Value = 1.285361611317384
$ws.Range("H6").Value = 1.004160822486213
$ws.Range("B7").Value = 0.1707995201049219
$ws.Range("C7").Value = 0.1428564451682919
$ws.Range("D7").Value = -0.9151433756610738
$ws.Range("E7").Value = 0.2330363382688319
$ws.Range("F7").Value = -0.4369475057286532
$ws.Range("G7").Value = -0.7181482945598243
$ws.Range("B8").Value = -0.02794307493663006
$ws.Range("C8").Value = -1.085942895765996
$ws.Range("D8").Value = 0.06223681816390991
$ws.Range("E8").Value = -0.6077470258335751
$ws.Range("F8").Value = -0.8889478146647463
$ws.Range("G8").Value = -1.028119891104211
$ws.Range("H8").Value = -1.248197453708784
$ws.Range("I8").Value = -1.109183561972742
$ws.Range("B9").Value = -1.057999820829366
$ws.Range("C9").Value = 0.09017989310053998
$ws.Range("D9").Value = -0.579803950896945
$ws.Range("E9").Value = -0.8610047397281162
$ws.Range("F9").Value = -1.000176816167581
$ws.Range("G9").Value = -1.220254378772154
$ws.Range("H9").Value = -1.081240487036112
$ws.Range("B10").Value = 1.148179713929906
$ws.Range("C10").Value = 0.4781958699324207
$ws.Range("D10").Value = 0.1969950811012495
$ws.Range("E10").Value = 0.05782300466178469
$ws.Range("F10").Value = -0.1622545579427879
$ws.Range("G10").Value = -0.02324066620674611
$ws.Range("B11").Value = -0.669983843997485
$ws.Range("C11").Value = -0.9511846328286562
$ws.Range("D11").Value = -1.090356709268121
$ws.Range("E11").Value = -1.310434271872694
$ws.Range("F11").Value = -1.171420380136652
$ws.Range("B12").Value = -0.2812007888311712
$ws.Range("C12").Value = -0.420372865270636
$ws.Range("D12").Value = -0.6404504278752086
$ws.Range("E12").Value = -0.5014365361391668
$ws.Range("B13").Value = -0.1391720764394648
$ws.Range("C13").Value = -0.3592496390440373
$ws.Range("D13").Value = -0.2202357473079956
$ws.Range("B14").Value = -0.2200775626045726
$ws.Range("C14").Value = -0.08106367086853081
$ws.Range("B15").Value = 0.1390138917360418
